$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data range (A1:A31) below the header, then write the new consolidated rows
$ws.Range("A2:A31").ClearContents()

$values = @(
    "('Construct', ['Token Artifact Creature — Construct', 'Defender', '1/1'])",
    "('Dack Fayden Emblem', ['Emblem — Dack', 'Whenever you cast a spell that targets one or more permanents, gain control of those permanents.'])",
    "('Demon', ['Token Creature — Demon', 'Flying', '*/*'])",
    "('Elephant', ['Token Creature — Elephant', '3/3'])",
    "('Ogre', ['Token Creature — Ogre', '4/4'])",
    "('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])",
    "('Squirrel', ['Token Creature — Squirrel', '1/1'])",
    "('Wolf', ['Token Creature — Wolf', '2/2'])",
    "('Zombie', ['Token Creature — Zombie', '2/2'])"
)

$row = 2
foreach ($val in $values) {
    $ws.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}
